$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 165 (shifts rows 165..259 down to 166..260)
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A165").Value = 10
$ws.Range("B165").Value = "Vega Modelo de Temuco"
$ws.Range("C165").Value = "La Araucanía"
$ws.Range("D165").Value = 44582
$ws.Range("E165").Value = 9
$ws.Range("F165").Value = 100112009
$ws.Range("G165").Value = "Acelga"
$ws.Range("H165").Value = "Sin especificar"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 50
$ws.Range("K165").Value = 7000
$ws.Range("L165").Value = 7000
$ws.Range("M165").Value = 7000
$ws.Range("N165").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O165").Value = "Provincia de Cautín"
$ws.Range("P165").Value = 583
$ws.Range("Q165").Value = 12
$ws.Range("R165").Value = "Hortaliza"
